$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '60.984.86'
$ws.Range('E2').Value = '  -4.04%  '
$ws.Range('D3').Value = "'" + '2.960.83'
$ws.Range('E3').Value = '  -3.96%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'" + '542.25'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('D6').Value = "'" + '128.32'
$ws.Range('E6').Value = '  -7.81%  '
$ws.Range('D8').Value = "'" + '2.953.39'
$ws.Range('E8').Value = '  -4.01%  '
$ws.Range('D9').Value = "'" + '0.487'
$ws.Range('E9').Value = '  -2.68%  '
$ws.Range('E10').Value = '  -7.76%  '
$ws.Range('D11').Value = "'" + '0.141'
$ws.Range('E11').Value = '  -8.91%  '
$ws.Range('D12').Value = "'" + '0.436'
$ws.Range('E12').Value = '  -4.75%  '
$ws.Range('E13').Value = '  -4.45%  '
$ws.Range('D14').Value = "'" + '33.06'
$ws.Range('E14').Value = '  -5.74%  '
$ws.Range('D15').Value = "'" + '3.445.03'
$ws.Range('E15').Value = '  -3.76%  '
$ws.Range('D16').Value = "'" + '61.045.09'
$ws.Range('E16').Value = '  -3.83%  '
$ws.Range('E17').Value = '  -3.21%  '
$ws.Range('D18').Value = "'" + '2.963.53'
$ws.Range('E18').Value = '  -3.83%  '
$ws.Range('D19').Value = "'" + '6.49'
$ws.Range('E19').Value = '  -2.64%  '
$ws.Range('D20').Value = "'" + '469.20'
$ws.Range('E20').Value = '  -1.38%  '
$ws.Range('D21').Value = "'" + '12.86'
$ws.Range('E21').Value = '  -4.73%  '
$ws.Range('D22').Value = "'" + '0.651'
$ws.Range('E22').Value = '  -7.26%  '
$ws.Range('D23').Value = "'" + '6.84'
$ws.Range('E23').Value = '  -3.70%  '
$ws.Range('D24').Value = "'" + '78.94'
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').Value = "'" + '11.77'
$ws.Range('E25').Value = '  -4.02%  '
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('D27').Value = "'" + '2.66'
$ws.Range('E27').Value = '  -2.39%  '
$ws.Range('D28').Value = "'" + '7.48'
$ws.Range('E28').Value = '  -6.24%  '
$ws.Range('D29').Value = "'" + '1.00'
$ws.Range('E29').Value = '  +0.30%  '
$ws.Range('D30').Value = "'" + '1.85'
$ws.Range('E30').Value = '  -2.73%  '
$ws.Range('D31').Value = "'" + '25.08'
$ws.Range('E31').Value = '  -4.46%  '
$ws.Range('D32').Value = "'" + '1.11'
$ws.Range('E32').Value = '  -4.30%  '
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').Value = "'" + '54.26'
$ws.Range('E33').Value = '  -6.50%  '
$ws.Range('B34').Value = 'Stacks'
$ws.Range('C34').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D34').Value = "'" + '2.24'
$ws.Range('E34').Value = '  -3.56%  '
$ws.Range('D35').Value = "'" + '5.31'
$ws.Range('E35').Value = '  -3.03%  '
$ws.Range('D36').Value = "'" + '5.75'
$ws.Range('E36').Value = '  -4.60%  '
$ws.Range('D37').Value = "'" + '440.38'
$ws.Range('E37').Value = '  -10.98%  '
$ws.Range('D38').Value = "'" + '3.080.62'
$ws.Range('E38').Value = '  -5.61%  '
$ws.Range('D39').Value = "'" + '0.0774'
$ws.Range('E39').Value = '  -3.32%  '
$ws.Range('E40').Value = '  -8.27%  '
$ws.Range('E41').Value = '  -3.85%  '
$ws.Range('D42').Value = "'" + '7.91'
$ws.Range('E42').Value = '  -3.02%  '
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('D44').Value = "'" + '2.27'
$ws.Range('E44').Value = '  -13.45%  '
$ws.Range('D45').Value = "'" + '24.89'
$ws.Range('E45').Value = '  -2.34%  '
$ws.Range('D46').Value = "'" + '0.235'
$ws.Range('E46').Value = '  -7.84%  '
$ws.Range('E47').Value = '  -3.34%  '
$ws.Range('D48').Value = "'" + '1.30'
$ws.Range('E48').Value = '  +8.82%  '
$ws.Range('E49').Value = '  -7.39%  '
$ws.Range('D50').Value = "'" + '113.17'
$ws.Range('E50').Value = '  -8.46%  '
$ws.Range('D51').Value = "'" + '0.0₃0468'
$ws.Range('E51').Value = '  -12.18%  '
